$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 208 (existing rows 208-229 shift down to 209-230)
$ws.Rows.Item(208).Insert()

# Populate the newly inserted row 208 with the new record
$ws.Cells.Item(208, 1).Value = 11
$ws.Cells.Item(208, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(208, 3).Value = "Bíobío"
$ws.Cells.Item(208, 4).Value = 45132
$ws.Cells.Item(208, 5).Value = 8
$ws.Cells.Item(208, 6).Value = 100112043
$ws.Cells.Item(208, 7).Value = "Pepino ensalada"
$ws.Cells.Item(208, 8).Value = "Sin especificar"
$ws.Cells.Item(208, 9).Value = "Primera"
$ws.Cells.Item(208, 10).Value = 100
$ws.Cells.Item(208, 11).Value = 10000
$ws.Cells.Item(208, 12).Value = 11000
$ws.Cells.Item(208, 13).Value = 10500
$ws.Cells.Item(208, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(208, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(208, 16).Value = 175
$ws.Cells.Item(208, 17).Value = 60
$ws.Cells.Item(208, 18).Value = "Hortaliza"
